$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.103.16'
$ws.Range("E2").Value = '  -3.09%  '
$ws.Range("D3").Value = '3.509.98'
$ws.Range("E3").Value = '  -5.10%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '577.56'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = '171.05'
$ws.Range("E6").Value = '  -3.88%  '
$ws.Range("D7").Value = '3.502.21'
$ws.Range("E7").Value = '  -5.03%  '
$ws.Range("D8").Value = '0.607'
$ws.Range("E8").Value = '  -1.41%  '
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("E10").Value = '  -5.75%  '
$ws.Range("D11").Value = '6.77'
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("E12").Value = '  -4.62%  '
$ws.Range("D13").Value = '47.12'
$ws.Range("E13").Value = '  -4.09%  '
$ws.Range("D14").Value = '0.0000272'
$ws.Range("E14").Value = '  -5.45%  '
$ws.Range("D15").Value = '4.080.33'
$ws.Range("E15").Value = '  -5.02%  '
$ws.Range("E16").Value = '  -5.64%  '
$ws.Range("D17").Value = '623.21'
$ws.Range("E17").Value = '  -8.18%  '
$ws.Range("D18").Value = '3.516.58'
$ws.Range("E18").Value = '  -4.92%  '
$ws.Range("D19").Value = '69.076.91'
$ws.Range("E19").Value = '  -3.32%  '
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("D21").Value = '17.39'
$ws.Range("E21").Value = '  -3.18%  '
$ws.Range("D22").Value = '11.14'
$ws.Range("E22").Value = '  -3.91%  '
$ws.Range("E23").Value = '  -6.39%  '
$ws.Range("D24").Value = '15.87'
$ws.Range("E24").Value = '  -9.04%  '
$ws.Range("D25").Value = '97.40'
$ws.Range("E25").Value = '  -4.71%  '
$ws.Range("E26").Value = '  -4.84%  '
$ws.Range("E28").Value = '  -6.91%  '
$ws.Range("D29").Value = '9.31'
$ws.Range("E29").Value = '  -9.64%  '
$ws.Range("D30").Value = '32.57'
$ws.Range("E30").Value = '  -7.52%  '
$ws.Range("D31").Value = '3.15'
$ws.Range("E31").Value = '  -8.05%  '
$ws.Range("D32").Value = '8.52'
$ws.Range("E32").Value = '  -7.28%  '
$ws.Range("E33").Value = '  -7.57%  '
$ws.Range("D34").Value = '7.00'
$ws.Range("E34").Value = '  -7.71%  '
$ws.Range("D35").Value = '634.38'
$ws.Range("E35").Value = '  +8.30%  '
$ws.Range("D36").Value = '10.72'
$ws.Range("E36").Value = '  -4.22%  '
$ws.Range("E37").Value = '  -5.69%  '
$ws.Range("D38").Value = '3.43'
$ws.Range("E38").Value = '  -16.23%  '
$ws.Range("E39").Value = '  -3.73%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").Value = '0.0448'
$ws.Range("E41").Value = '  -2.25%  '
$ws.Range("E42").Value = '  -6.45%  '
$ws.Range("D43").Value = '3.380.64'
$ws.Range("E43").Value = '  -8.22%  '
$ws.Range("D44").Value = '0.326'
$ws.Range("E44").Value = '  -7.35%  '
$ws.Range("D45").Value = '32.83'
$ws.Range("E45").Value = '  -8.09%  '
$ws.Range("D46").Value = '0.0₃0689'
$ws.Range("E46").Value = '  -10.26%  '
$ws.Range("E47").Value = '  -7.44%  '
$ws.Range("E48").Value = '  -5.51%  '
$ws.Range("E49").Value = '  -2.76%  '
$ws.Range("E50").Value = '  +14.47%  '
$ws.Range("D51").Value = '131.94'
$ws.Range("E51").Value = '  -2.51%  '
